$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 7 through 26 (the rows that get collapsed into the grouped tuples)
$ws.Range("A7:A26").EntireRow.Delete() | Out-Null

# Update the remaining rows (A2:A6) with the new grouped tuple-style strings
$ws.Range("A2").Value = "('Eldrazi Spawn', ['Token Creature — Eldrazi Spawn', 'Sacrifice this creature: Add {C}.', '0/1'])"
$ws.Range("A3").Value = "('Elemental', ['Token Creature — Elemental', '*/*'])"
$ws.Range("A4").Value = "('Hellion', ['Token Creature — Hellion', '4/4'])"
$ws.Range("A5").Value = "('Ooze', ['Token Creature — Ooze', '*/*'])"
$ws.Range("A6").Value = "('Tuktuk the Returned', ['Token Artifact Creature — Goblin Golem', 'Tuktuk the Returned is legendary.', '5/5'])"
